# Switch to newest InOutModule format:
# set an explicit custom height (24pt) on the header row (row 1)
# for every worksheet in the workbook.
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Rows.Item(1).RowHeight = 24
}
